# Apply updates to the "Work Report" worksheet to refresh the generated
# report with the latest billing totals for the single work request
# contained in this file, and to clear the (no longer applicable) Scope ID.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Report Generated On" timestamp (D5)
$ws.Range("D5").Value = "Report Generated On: 08/26/2025 10:02 AM"

# Update the Total Billed Amount summary value (C8)
$ws.Range("C8").Value = 478.55

# Clear the Scope ID # value (G10) - no longer populated
$ws.Range("G10").Value = ""

# Update the line item pricing (H16) and the TOTAL row (H17)
$ws.Range("H16").Value = 478.55
$ws.Range("H17").Value = 478.55
